# Insert a new weekly price record as row 228 (shifting all subsequent
# data rows down by one, 228-299 -> 229-300), matching the "Fruta /
# hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 228 (and everything below it) down by one row.
$ws.Rows("228:228").Insert()

# Populate the newly inserted row 228 with the new weekly record.
$ws.Range("A228").Value = 10
$ws.Range("B228").Value = "Vega Modelo de Temuco"
$ws.Range("C228").Value = "La Araucanía"
$ws.Range("D228").Value = 45146
$ws.Range("E228").Value = 9
$ws.Range("F228").Value = 100112012
$ws.Range("G228").Value = "Espinaca"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 80
$ws.Range("K228").Value = 8000
$ws.Range("L228").Value = 8000
$ws.Range("M228").Value = 8000
$ws.Range("N228").Value = "$/docena de paquetes"
$ws.Range("O228").Value = "Región de La Araucanía"
$ws.Range("P228").Value = 667
$ws.Range("Q228").Value = 12
$ws.Range("R228").Value = "Hortaliza"
